# MODIFICACION DEL MODAL DE IMAGENES
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old "COMENTADO CON ERDNANADO" note is replaced with a short "it " note
$ws.Range("C37").Value = "it "

# Append a new pending task row at the bottom of the list
$ws.Range("A100").Value = "CAMBIAR ALERT POR MODAL"
$ws.Range("A100").WrapText = $true
$ws.Range("A100").HorizontalAlignment = -4131
$ws.Range("A100").VerticalAlignment = -4160

# Row 18 reflows to a shorter, auto-fit height now that the layout changed
$ws.Rows.Item(18).RowHeight = 18.75

# Scroll the view down and select the new last cell, matching the new bottom of the list
$excel.ActiveWindow.ScrollRow = 87
$ws.Range("B100").Select()
